# Add a new bullet point right after the bullet that ends in
# "...reducing UI-related bug issues by ~50%" and right before the bullet
# that starts with "Cut tabular data loading time...".
#
# New bullet text:
#   "Reduced CI pipeline cost by ~70% and decreased test suite execution
#    time from 26m to 5m by introducing parallelisation for Jest tests,
#    linting, formatting, and building steps on GitHub Actions"

$d = $word.ActiveDocument

# Find the anchor paragraph using a short, unique substring instead of a
# hard-coded paragraph index, so the script keeps working even if the
# surrounding content shifts slightly.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("UI-related bug issues by", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph ('UI-related bug issues by') to insert the new bullet after"
}

# Find.Execute collapses/positions $searchRange on the found text; the
# paragraph that contains it is the "...~50%" bullet we want to insert after.
$anchorPara = $searchRange.Paragraphs(1)

# Create a new (initially empty) paragraph right after the anchor paragraph.
# It automatically inherits the anchor's paragraph formatting (same
# ListParagraph style / numPr bullet list as its neighbours).
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()

# Replace the (empty) new paragraph's contents with the fully-formatted
# OOXML for the new bullet point, run by run, matching the target markup.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Reduced </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">CI pipeline cost by </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Inter Light" w:hAnsi="Inter Light" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>~</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>70% and decreas</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ed</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> test suite execution time from 26</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">m </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">to </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>5m</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> by introducing</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>parallelisation for Jest tests, lint</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ing</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>, formatting, and build</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ing</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> steps on GitHub Actions</w:t></w:r></w:p>
'@

[void]$newPara.Range.InsertXML($xml)
